$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column F (Dependency_Type_Descrip / "Business Description..") entirely.
# Remaining columns G (Dependency_Name) and H (Dependency_Descrip) shift left
# to become F and G respectively.
$ws.Columns("F").Delete()

# Select the (new) column F, matching the post-edit selection in the file.
$ws.Columns("F").Select() | Out-Null
